$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 1 header values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON) values
$ws.Range("B2").Value = 8.3808187752763956
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 11.471684398943353
$ws.Range("E2").ClearContents()

# Row 3 (STR) values
$ws.Range("B3").Value = 6.877254376254939
$ws.Range("C3").Value = -6.2925999639750998
$ws.Range("D3").Value = 11.964927595038091
$ws.Range("E3").Value = -5.548955410987837

# Update selection to match recorded view state
$ws.Range("B1:E3").Select()
